# Update the "取得日時" (acquired datetime) timestamp column for all data
# rows on the "ランサーズ" sheet from 2026-01-06 18:28:30 to 2026-01-06 18:37:18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-06 18:37:18"

# Find the last used row in column A and update rows 2..lastRow (row 1 is header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}
